$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "египетских фунтов"
$ws.Range("B2").Value = 26.113800048828125

$ws.Range("A3").Value = "фунт стерлингов соединенного королевства"
$ws.Range("B3").Value = 99.50350189208984
